# Applies the cryptocurrency price/volume refresh captured in the
# "Updated cryptos list on Sat Nov  2 14:37:57 UTC 2024 with GitHub Actions" commit.
# Cells D (Price) and E (Volume(1h)) are plain text in the workbook (not numbers),
# so values are written through a text-literal formula + paste-values round trip to
# avoid Excel auto-converting numeric-looking strings (e.g. "0.510", "1.00") into
# real numbers and to avoid touching each cell's existing style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues: keep text, drop formula/format churn
}

$excel.CutCopyMode = $false

# Row 2: update D2, E2
Set-TextValue $ws.Range("D2") '69.221.40'
Set-TextValue $ws.Range("E2") '  -3.17%  '

# Row 3: update D3, E3
Set-TextValue $ws.Range("D3") '2.482.63'
Set-TextValue $ws.Range("E3") '  -3.53%  '

# Row 4: update E4
Set-TextValue $ws.Range("E4") '  -0.01%  '

# Row 5: update D5, E5
Set-TextValue $ws.Range("D5") '565.58'
Set-TextValue $ws.Range("E5") '  -3.14%  '

# Row 6: update D6, E6
Set-TextValue $ws.Range("D6") '163.75'
Set-TextValue $ws.Range("E6") '  -5.31%  '

# Row 7: update E7
Set-TextValue $ws.Range("E7") '  -0.04%  '

# Row 8: update D8, E8
Set-TextValue $ws.Range("D8") '0.510'
Set-TextValue $ws.Range("E8") '  -1.83%  '

# Row 9: update D9, E9
Set-TextValue $ws.Range("D9") '2.483.09'
Set-TextValue $ws.Range("E9") '  -3.38%  '

# Row 10: update D10, E10
Set-TextValue $ws.Range("D10") '0.157'
Set-TextValue $ws.Range("E10") '  -5.99%  '

# Row 11: update E11
Set-TextValue $ws.Range("E11") '  -0.77%  '

# Row 12: update D12, E12
Set-TextValue $ws.Range("D12") '0.351'
Set-TextValue $ws.Range("E12") '  -2.71%  '

# Row 13: update D13, E13
Set-TextValue $ws.Range("D13") '4.89'
Set-TextValue $ws.Range("E13") '  -0.79%  '

# Row 14: update D14, E14
Set-TextValue $ws.Range("D14") '2.937.34'
Set-TextValue $ws.Range("E14") '  -3.57%  '

# Row 15: update D15, E15
Set-TextValue $ws.Range("D15") '69.169.05'
Set-TextValue $ws.Range("E15") '  -3.07%  '

# Row 16: update D16, E16
Set-TextValue $ws.Range("D16") '0.0000174'
Set-TextValue $ws.Range("E16") '  -3.88%  '

# Row 17: update D17, E17
Set-TextValue $ws.Range("D17") '24.16'
Set-TextValue $ws.Range("E17") '  -5.53%  '

# Row 18: update D18, E18
Set-TextValue $ws.Range("D18") '2.485.32'
Set-TextValue $ws.Range("E18") '  -3.67%  '

# Row 19: update D19, E19
Set-TextValue $ws.Range("D19") '11.11'
Set-TextValue $ws.Range("E19") '  -4.87%  '

# Row 20: update D20, E20
Set-TextValue $ws.Range("D20") '7.35'
Set-TextValue $ws.Range("E20") '  -7.73%  '

# Row 21: update D21, E21
Set-TextValue $ws.Range("D21") '345.04'
Set-TextValue $ws.Range("E21") '  -3.88%  '

# Row 22: update D22, E22
Set-TextValue $ws.Range("D22") '3.84'
Set-TextValue $ws.Range("E22") '  -3.46%  '

# Row 23: update D23
Set-TextValue $ws.Range("D23") '1.91'

# Row 25: update D25, E25
Set-TextValue $ws.Range("D25") '69.40'
Set-TextValue $ws.Range("E25") '  -1.89%  '

# Row 26: update D26, E26
Set-TextValue $ws.Range("D26") '3.86'
Set-TextValue $ws.Range("E26") '  -6.53%  '

# Row 27: update B27, C27, D27, E27
Set-TextValue $ws.Range("B27") 'WrappedeETH'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range("D27") '2.609.43'
Set-TextValue $ws.Range("E27") '  -3.60%  '

# Row 28: update B28, C28, D28, E28
Set-TextValue $ws.Range("B28") 'Aptos'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D28") '8.64'
Set-TextValue $ws.Range("E28") '  -6.12%  '

# Row 29: update D29, E29
Set-TextValue $ws.Range("D29") '0.998'
Set-TextValue $ws.Range("E29") '  -0.27%  '

# Row 30: update D30, E30
Set-TextValue $ws.Range("D30") '0.0₃0868'
Set-TextValue $ws.Range("E30") '  -6.58%  '

# Row 31: update D31, E31
Set-TextValue $ws.Range("D31") '7.65'
Set-TextValue $ws.Range("E31") '  -4.52%  '

# Row 32: update D32, E32
Set-TextValue $ws.Range("D32") '440.57'
Set-TextValue $ws.Range("E32") '  -7.83%  '

# Row 33: update D33, E33
Set-TextValue $ws.Range("D33") '1.18'
Set-TextValue $ws.Range("E33") '  -8.52%  '

# Row 34: update D34, E34
Set-TextValue $ws.Range("D34") '1.00'
Set-TextValue $ws.Range("E34") '  +0.06%  '

# Row 35: update E35
Set-TextValue $ws.Range("E35") '  -4.89%  '

# Row 36: update D36, E36
Set-TextValue $ws.Range("D36") '154.87'
Set-TextValue $ws.Range("E36") '  -2.05%  '

# Row 37: update E37
Set-TextValue $ws.Range("E37") '  -4.98%  '

# Row 38: update E38
Set-TextValue $ws.Range("E38") '  -0.59%  '

# Row 39: update D39, E39
Set-TextValue $ws.Range("D39") '18.06'
Set-TextValue $ws.Range("E39") '  -4.57%  '

# Row 40: update E40
Set-TextValue $ws.Range("E40") '  +0.00%  '

# Row 41: update D41, E41
Set-TextValue $ws.Range("D41") '0.312'
Set-TextValue $ws.Range("E41") '  -3.58%  '

# Row 42: update E42
Set-TextValue $ws.Range("E42") '  -7.24%  '

# Row 43: update E43
Set-TextValue $ws.Range("E43") '  -4.48%  '

# Row 44: update D44, E44
Set-TextValue $ws.Range("D44") '37.89'
Set-TextValue $ws.Range("E44") '  -2.29%  '

# Row 45: update D45, E45
Set-TextValue $ws.Range("D45") '2.15'
Set-TextValue $ws.Range("E45") '  -10.51%  '

# Row 46: update E46
Set-TextValue $ws.Range("E46") '  -9.81%  '

# Row 47: update D47, E47
Set-TextValue $ws.Range("D47") '138.53'
Set-TextValue $ws.Range("E47") '  -5.52%  '

# Row 48: update D48, E48
Set-TextValue $ws.Range("D48") '3.42'
Set-TextValue $ws.Range("E48") '  -4.27%  '

# Row 49: update D49, E49
Set-TextValue $ws.Range("D49") '0.510'
Set-TextValue $ws.Range("E49") '  -6.04%  '

# Row 50: update E50
Set-TextValue $ws.Range("E50") '  -2.72%  '

# Row 51: update D51, E51
Set-TextValue $ws.Range("D51") '0.571'
Set-TextValue $ws.Range("E51") '  -3.11%  '
